# Weekly fruit/vegetable price update.
# A new price record (row) is inserted at row 36, pushing the existing
# records (previously rows 36-64) down to rows 37-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36 (shifts rows 36:64 -> 37:65)
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new weekly record
$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value = "Bíobío"
$ws.Cells.Item(36, 4).Value = 44484
$ws.Cells.Item(36, 5).Value = 8
$ws.Cells.Item(36, 6).Value = 100112032
$ws.Cells.Item(36, 7).Value = "Zapallo italiano"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 450
$ws.Cells.Item(36, 11).Value = 12000
$ws.Cells.Item(36, 12).Value = 13000
$ws.Cells.Item(36, 13).Value = 12556
$ws.Cells.Item(36, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(36, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(36, 16).Value = 251
$ws.Cells.Item(36, 17).Value = 50
$ws.Cells.Item(36, 18).Value = "Hortaliza"
